$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates (rows 3-5: E/F columns) ---
$ws = $wb.Worksheets.Item("Schedule")
$ws.Range("E3").Value = 52.26773175000002
$ws.Range("F3").Value = 1.455520238095239
$ws.Range("E4").Value = 471.5901352499999
$ws.Range("F4").Value = 31.18982375992063
$ws.Range("E5").Value = -60.296652
$ws.Range("F5").Value = -1.77238835978836

# --- Detailed sheet updates (rows 32-97: B/C columns) ---
$ws = $wb.Worksheets.Item("Detailed")
$ws.Range("B32").Value = 0.00016
$ws.Range("B33").Value = 0.009429999999999999
$ws.Range("B34").Value = 0.51
$ws.Range("C34").Value = "historical"
$ws.Range("B35").Value = 0.7
$ws.Range("B36").Value = 4.80008
$ws.Range("B37").Value = 8.264950000000001
$ws.Range("B38").Value = 26.89449
$ws.Range("B39").Value = 57.19078
$ws.Range("B40").Value = 60.89812
$ws.Range("B41").Value = 64.93597
$ws.Range("B42").Value = 68.31175
$ws.Range("B43").Value = 71.46543
$ws.Range("B44").Value = 67.32592
$ws.Range("B45").Value = 66.65888
$ws.Range("B46").Value = 64.8901
$ws.Range("B47").Value = 59.37938
$ws.Range("B48").Value = 60.76462
$ws.Range("B49").Value = 64.04767
$ws.Range("B50").Value = 63.00946
$ws.Range("B51").Value = 61.57869
$ws.Range("B52").Value = 63.10166
$ws.Range("B59").Value = 57.49046
$ws.Range("B60").Value = 65
$ws.Range("B61").Value = 77.94
$ws.Range("B62").Value = 78
$ws.Range("B63").Value = 66.67346999999999
$ws.Range("B65").Value = 8.718500000000001
$ws.Range("B66").Value = 0.5099399999999999
$ws.Range("B68").Value = -2.55434
$ws.Range("B69").Value = -5.62917
$ws.Range("B70").Value = -6.96628
$ws.Range("B71").Value = -6.60898
$ws.Range("B72").Value = -8.66264
$ws.Range("B73").Value = -8.643409999999999
$ws.Range("B74").Value = -9.68276
$ws.Range("B75").Value = -10
$ws.Range("B76").Value = -9.551170000000001
$ws.Range("B77").Value = -8.70303
$ws.Range("B78").Value = -8.271380000000001
$ws.Range("B79").Value = -7.88906
$ws.Range("B80").Value = -7.96894
$ws.Range("B81").Value = -6
$ws.Range("B82").Value = -2.84369
$ws.Range("B83").Value = -5.5912
$ws.Range("B84").Value = -2.96533
$ws.Range("B85").Value = 0.01012
$ws.Range("B86").Value = 9.522930000000001
$ws.Range("B87").Value = 46.32246
$ws.Range("B88").Value = 57.3
$ws.Range("B89").Value = 69.61881
$ws.Range("B90").Value = 65
$ws.Range("B92").Value = 59.15624
$ws.Range("B93").Value = 58.33176
$ws.Range("B95").Value = 59.93754
$ws.Range("B96").Value = 63.27623
$ws.Range("B97").Value = 64.12067

